$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Category" (B1) to "Vendor"
$ws.Range("B1").Value = "Vendor"

# Remove the now-duplicate "Vendor" header that was in I1
$ws.Range("I1").ClearContents()
